# Re-process the metadata sheet with the newly curated dimensions/measures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "superficie-energias-renovables-en-retirada" becomes a measure instead of a dimension
$ws.Range("B2").Value = "iaest-measure:superficie-energias-renovables-en-retirada"
$ws.Range("B3").Value = "medida"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("B5").Clear()

# "municipio-nombre" becomes a dimension (refArea) instead of a measure
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# "otros-cultivos-para-renovables" becomes a measure instead of a dimension
$ws.Range("F2").Value = "iaest-measure:otros-cultivos-para-renovables"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("F5").Clear()

Write-Output "edit applied"
